# Apply updated pl_mw values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.7213517308646544
    "C2" = 0.1718431795037976
    "D2" = 0.06496477394830613
    "F2" = 1.511035748728276
    "G2" = 1.41764259692188
    "H2" = 1.288948255308114
    "K2" = 0.3301834923238971
    "L2" = 0.3050637716089284
    "B3" = 0.6824973833680588
    "C3" = 0.1717174886374124
    "D3" = 0.06474175850884123
    "F3" = 1.494950897819109
    "G3" = 1.401719831964201
    "H3" = 1.286134612341229
    "K3" = 0.2983940222513581
    "L3" = 0.2942318151794296
    "B4" = 0.6590565873403023
    "C4" = 0.1716358574595223
    "D4" = 0.06459688050109236
    "F4" = 1.485791888115344
    "G4" = 1.392627433488855
    "H4" = 1.284895940418238
    "K4" = 0.2790127680990082
    "L4" = 0.2877618200802203
    "B5" = 0.6496090867790372
    "C5" = 0.1716014723689128
    "D5" = 0.06453584962459047
    "F5" = 1.482239630587856
    "G5" = 1.389093932265439
    "H5" = 1.284514069958973
    "K5" = 0.2711494874784961
    "L5" = 0.2851706721039022
    "B6" = 0.6480466738018436
    "C6" = 0.1715956951510123
    "D6" = 0.0645255953687176
    "F6" = 1.48166065290998
    "G6" = 1.388517559112088
    "H6" = 1.284458082175703
    "K6" = 0.2698458978896809
    "L6" = 0.2847431561383047
    "B7" = 0.6589287503522314
    "C7" = 0.1716353982624117
    "D7" = 0.06459606547575447
    "F7" = 1.485743252154791
    "G7" = 1.392579084595155
    "H7" = 1.284890292832614
    "K7" = 0.2789065803118405
    "L7" = 0.2877266910867462
    "B8" = 0.7078685450143212
    "C8" = 0.1718007667288255
    "D8" = 0.06488953089796823
    "F8" = 1.505340705730603
    "G8" = 1.412010204318236
    "H8" = 1.287876606923348
    "K8" = 0.319193991686177
    "L8" = 0.3012913368137191
    "B9" = 0.8071357546114086
    "C9" = 0.172089663513276
    "D9" = 0.06540172063954586
    "F9" = 1.549475126859335
    "G9" = 1.455562197576796
    "H9" = 1.297615617899424
    "K9" = 0.3992879050318265
    "L9" = 0.3293309133720186
    "B10" = 0.8820808639974018
    "C10" = 0.1722803140321894
    "D10" = 0.06573912850587149
    "F10" = 1.585401943901203
    "G10" = 1.490911343895732
    "H10" = 1.307145106672891
    "K10" = 0.4588023491721174
    "L10" = 0.3508175090868377
    "B11" = 0.9166140273288761
    "C11" = 0.1723623516588191
    "D11" = 0.06588411613811473
    "F11" = 1.602511743926357
    "G11" = 1.50772718989424
    "H11" = 1.311997558567356
    "K11" = 0.4860238500847629
    "L11" = 0.3607866767967067
    "B12" = 0.9297540808586291
    "C12" = 0.1723927423136722
    "D12" = 0.06593779154794177
    "F12" = 1.60910133953405
    "G12" = 1.514201105653001
    "H12" = 1.313909559136505
    "K12" = 0.4963531996169763
    "L12" = 0.364589860963747
    "B13" = 0.9269213324603811
    "C13" = 0.1723862271839511
    "D13" = 0.06592628630538044
    "F13" = 1.607677234156242
    "G13" = 1.512802106324699
    "H13" = 1.313494461774752
    "K13" = 0.49412765050036
    "L13" = 0.3637695274645978
    "B14" = 0.9176938041528047
    "C14" = 0.1723648654529324
    "D14" = 0.0658885566950218
    "F14" = 1.60305165812261
    "G14" = 1.508257673873004
    "H14" = 1.312153366835446
    "K14" = 0.4868732286169006
    "L14" = 0.361099003980101
    "B15" = 0.9120498870541951
    "C15" = 0.1723516928217208
    "D15" = 0.06586528609339481
    "F15" = 1.60023275552453
    "G15" = 1.50548790894544
    "H15" = 1.311341609503671
    "K15" = 0.4824324371597299
    "L15" = 0.3594668895562734
    "B16" = 0.8798328812167711
    "C16" = 0.1722748582333402
    "D16" = 0.06572948165979042
    "F16" = 1.584299224811247
    "G16" = 1.48982721813536
    "H16" = 1.306838406146198
    "K16" = 0.457026328133793
    "L16" = 0.3501699259773261
    "B17" = 0.860181397927505
    "C17" = 0.1722265208460954
    "D17" = 0.06564398853166864
    "F17" = 1.574721028250622
    "G17" = 1.480408499696324
    "H17" = 1.30420841871782
    "K17" = 0.4414783324118048
    "L17" = 0.3445164798380773
    "B18" = 0.8489198360947228
    "C18" = 0.1721982770291834
    "D18" = 0.06559401533500164
    "F18" = 1.569284034836272
    "G18" = 1.475060293204251
    "H18" = 1.302744418819401
    "K18" = 0.4325494825533553
    "L18" = 0.3412830995061995
    "B19" = 0.845113990724343
    "C19" = 0.1721886383675209
    "D19" = 0.06557695808271191
    "F19" = 1.567455543027307
    "G19" = 1.473261354512857
    "H19" = 1.302257095492052
    "K19" = 0.429528726865442
    "L19" = 0.3401914791567862
    "B20" = 0.8622690433578271
    "C20" = 0.1722317121349271
    "D20" = 0.06565317224067968
    "F20" = 1.575733176337877
    "G20" = 1.481403974573027
    "H20" = 1.304483344624288
    "K20" = 0.4431320007347779
    "L20" = 0.345116401043839
    "B21" = 0.9204024422251678
    "C21" = 0.1723711582413934
    "D21" = 0.06589967217615111
    "F21" = 1.604407301288177
    "G21" = 1.509589600764144
    "H21" = 1.312545256895135
    "K21" = 0.489003455019315
    "L21" = 0.3618826389564873
    "B22" = 0.9587636434100659
    "C22" = 0.172458358384457
    "D22" = 0.06605361233060592
    "F22" = 1.623791616449267
    "G22" = 1.528629260068698
    "H22" = 1.318248358133729
    "K22" = 0.5191063662878435
    "L22" = 0.3730040273859174
    "B23" = 0.9382559765908809
    "C23" = 0.1724121783944028
    "D23" = 0.06597210875595394
    "F23" = 1.613386819177407
    "G23" = 1.518410697300908
    "H23" = 1.315164753569888
    "K23" = 0.503028641278064
    "L23" = 0.3670533348437601
    "B24" = 0.8613251056568458
    "C24" = 0.1722293665670449
    "D24" = 0.06564902284582885
    "F24" = 1.575275367276475
    "G24" = 1.480953712511308
    "H24" = 1.304358901052552
    "K24" = 0.4423843464418837
    "L24" = 0.344845124154844
    "B25" = 0.7799280555617543
    "C25" = 0.1720153012040821
    "D25" = 0.06526996892149839
    "F25" = 1.536922474914959
    "G25" = 1.44319387280504
    "H25" = 1.294564535048465
    "K25" = 0.3775031716853334
    "L25" = 0.3215905022711354
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
